$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.05973035074357
$ws.Cells.Item(2, 4).Value = 1.059168720306137
$ws.Cells.Item(2, 5).Value = 1.063424380384006
$ws.Cells.Item(2, 6).Value = 1.065024848663887
$ws.Cells.Item(2, 9).Value = 1.04967398027987
$ws.Cells.Item(2, 10).Value = 1.064714865573094
$ws.Cells.Item(2, 11).Value = 1.06189895930655
$ws.Cells.Item(2, 12).Value = 1.066143056990225
$ws.Cells.Item(2, 13).Value = 1.067739202757813
$ws.Cells.Item(2, 14).Value = 1.06622688280671
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.061256273571145
$ws.Cells.Item(3, 4).Value = 1.060384423101575
$ws.Cells.Item(3, 5).Value = 1.064901751393008
$ws.Cells.Item(3, 6).Value = 1.066701514120702
$ws.Cells.Item(3, 9).Value = 1.050222517550651
$ws.Cells.Item(3, 10).Value = 1.065891402775624
$ws.Cells.Item(3, 11).Value = 1.062927822446664
$ws.Cells.Item(3, 12).Value = 1.067433774530817
$ws.Cells.Item(3, 13).Value = 1.069229033734981
$ws.Cells.Item(3, 14).Value = 1.067405090827018
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.062241526993218
$ws.Cells.Item(4, 4).Value = 1.061169064399718
$ws.Cells.Item(4, 5).Value = 1.065855938953291
$ws.Cells.Item(4, 6).Value = 1.067784646492177
$ws.Cells.Item(4, 9).Value = 1.050574947199319
$ws.Cells.Item(4, 10).Value = 1.066650174246021
$ws.Cells.Item(4, 11).Value = 1.063591002257221
$ws.Cells.Item(4, 12).Value = 1.068266674012808
$ws.Cells.Item(4, 13).Value = 1.07019080212368
$ws.Cells.Item(4, 14).Value = 1.068164939839935
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.06265522933692
$ws.Cells.Item(5, 4).Value = 1.061498455698135
$ws.Cells.Item(5, 5).Value = 1.066256664132263
$ws.Cells.Item(5, 6).Value = 1.068239579213521
$ws.Cells.Item(5, 9).Value = 1.050722511425787
$ws.Cells.Item(5, 10).Value = 1.06696856425197
$ws.Cells.Item(5, 11).Value = 1.063869195788671
$ws.Cells.Item(5, 12).Value = 1.068616286714859
$ws.Cells.Item(5, 13).Value = 1.070594600080268
$ws.Cells.Item(5, 14).Value = 1.068483781996213
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.062724662841297
$ws.Cells.Item(6, 4).Value = 1.061553734471395
$ws.Cells.Item(6, 5).Value = 1.066323923548292
$ws.Cells.Item(6, 6).Value = 1.068315940322796
$ws.Cells.Item(6, 9).Value = 1.050747253200249
$ws.Cells.Item(6, 10).Value = 1.067021988509945
$ws.Cells.Item(6, 11).Value = 1.063915870282636
$ws.Cells.Item(6, 12).Value = 1.068674956890413
$ws.Cells.Item(6, 13).Value = 1.070662368763973
$ws.Cells.Item(6, 14).Value = 1.068537282122762
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.062247056850014
$ws.Cells.Item(7, 4).Value = 1.061173467588201
$ws.Cells.Item(7, 5).Value = 1.065861295083936
$ws.Cells.Item(7, 6).Value = 1.067790726947519
$ws.Cells.Item(7, 9).Value = 1.050576921301567
$ws.Cells.Item(7, 10).Value = 1.066654430927592
$ws.Cells.Item(7, 11).Value = 1.063594721870998
$ws.Cells.Item(7, 12).Value = 1.06827134766007
$ws.Cells.Item(7, 13).Value = 1.070196199760324
$ws.Cells.Item(7, 14).Value = 1.068169202566482
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.06024648834604
$ws.Cells.Item(8, 4).Value = 1.05957999093484
$ws.Cells.Item(8, 5).Value = 1.063924036534498
$ws.Cells.Item(8, 6).Value = 1.065591861015449
$ws.Cells.Item(8, 9).Value = 1.049859883193696
$ws.Cells.Item(8, 10).Value = 1.065113009480195
$ws.Cells.Item(8, 11).Value = 1.062247203084005
$ws.Cells.Item(8, 12).Value = 1.066579738081091
$ws.Cells.Item(8, 13).Value = 1.068243170192394
$ws.Cells.Item(8, 14).Value = 1.066625592123827
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.056704561698057
$ws.Cells.Item(9, 4).Value = 1.056756458526215
$ws.Cells.Item(9, 5).Value = 1.06049639715258
$ws.Cells.Item(9, 6).Value = 1.061703053243076
$ws.Cells.Item(9, 9).Value = 1.048576968169279
$ws.Cells.Item(9, 10).Value = 1.062377142989791
$ws.Cells.Item(9, 11).Value = 1.059852780258857
$ws.Cells.Item(9, 12).Value = 1.063581082027246
$ws.Cells.Item(9, 13).Value = 1.064784002695362
$ws.Cells.Item(9, 14).Value = 1.063885840389197
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.054331469523272
$ws.Cells.Item(10, 4).Value = 1.054863166677429
$ws.Cells.Item(10, 5).Value = 1.058201390388278
$ws.Cells.Item(10, 6).Value = 1.059100351520699
$ws.Cells.Item(10, 9).Value = 1.047708401880986
$ws.Cells.Item(10, 10).Value = 1.06053952891478
$ws.Cells.Item(10, 11).Value = 1.05824269604313
$ws.Cells.Item(10, 12).Value = 1.061569497286902
$ws.Cells.Item(10, 13).Value = 1.062465395955259
$ws.Cells.Item(10, 14).Value = 1.062045616691423
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.053300959481014
$ws.Cells.Item(11, 4).Value = 1.054040659164798
$ws.Cells.Item(11, 5).Value = 1.057205154785645
$ws.Cells.Item(11, 6).Value = 1.057970793414922
$ws.Cells.Item(11, 9).Value = 1.047329096021486
$ws.Cells.Item(11, 10).Value = 1.059740466652838
$ws.Cells.Item(11, 11).Value = 1.057542146652254
$ws.Cells.Item(11, 12).Value = 1.06069538728246
$ws.Cells.Item(11, 13).Value = 1.061458313959382
$ws.Cells.Item(11, 14).Value = 1.061245419669416
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.052917727987501
$ws.Cells.Item(12, 4).Value = 1.05373472919527
$ws.Cells.Item(12, 5).Value = 1.056834725047815
$ws.Cells.Item(12, 6).Value = 1.057550826054682
$ws.Cells.Item(12, 9).Value = 1.047187717722346
$ws.Cells.Item(12, 10).Value = 1.059443144930911
$ws.Cells.Item(12, 11).Value = 1.057281416796302
$ws.Cells.Item(12, 12).Value = 1.060370231405318
$ws.Cells.Item(12, 13).Value = 1.061083759824712
$ws.Cells.Item(12, 14).Value = 1.060947675716541
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.052999953184664
$ws.Cells.Item(13, 4).Value = 1.053800371085979
$ws.Cells.Item(13, 5).Value = 1.056914201006909
$ws.Cells.Item(13, 6).Value = 1.05764092875159
$ws.Cells.Item(13, 9).Value = 1.047218065972597
$ws.Cells.Item(13, 10).Value = 1.059506944875616
$ws.Cells.Item(13, 11).Value = 1.057337367651885
$ws.Cells.Item(13, 12).Value = 1.060439999982896
$ws.Cells.Item(13, 13).Value = 1.061164124826294
$ws.Cells.Item(13, 14).Value = 1.061011566264484
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.053269290779158
$ws.Cells.Item(14, 4).Value = 1.054015379393981
$ws.Cells.Item(14, 5).Value = 1.057174542831903
$ws.Cells.Item(14, 6).Value = 1.057936087000284
$ws.Cells.Item(14, 9).Value = 1.047317419614181
$ws.Cells.Item(14, 10).Value = 1.059715900493678
$ws.Cells.Item(14, 11).Value = 1.057520605185495
$ws.Cells.Item(14, 12).Value = 1.060668519491463
$ws.Cells.Item(14, 13).Value = 1.061427363038711
$ws.Cells.Item(14, 14).Value = 1.061220818623491
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.053435178026563
$ws.Cells.Item(15, 4).Value = 1.054147797987478
$ws.Cells.Item(15, 5).Value = 1.057334896892201
$ws.Cells.Item(15, 6).Value = 1.058117890454525
$ws.Cells.Item(15, 9).Value = 1.047378569966576
$ws.Cells.Item(15, 10).Value = 1.059844576506181
$ws.Cells.Item(15, 11).Value = 1.05763343544566
$ws.Cells.Item(15, 12).Value = 1.060809255000428
$ws.Cells.Item(15, 13).Value = 1.061589488956854
$ws.Cells.Item(15, 14).Value = 1.061349677370692
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.05439979803063
$ws.Cells.Item(16, 4).Value = 1.054917696192956
$ws.Cells.Item(16, 5).Value = 1.058267454065652
$ws.Cells.Item(16, 6).Value = 1.05917526130299
$ws.Cells.Item(16, 9).Value = 1.047733507101708
$ws.Cells.Item(16, 10).Value = 1.060592488408415
$ws.Cells.Item(16, 11).Value = 1.058289117469916
$ws.Cells.Item(16, 12).Value = 1.061627443361006
$ws.Cells.Item(16, 13).Value = 1.062532166141986
$ws.Cells.Item(16, 14).Value = 1.062098651393614
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.055004081891487
$ws.Cells.Item(17, 4).Value = 1.055399904272894
$ws.Cells.Item(17, 5).Value = 1.058851750512421
$ws.Cells.Item(17, 6).Value = 1.059837824670491
$ws.Cells.Item(17, 9).Value = 1.047955286752287
$ws.Cells.Item(17, 10).Value = 1.061060726930274
$ws.Cells.Item(17, 11).Value = 1.058699501016777
$ws.Cells.Item(17, 12).Value = 1.062139839745842
$ws.Cells.Item(17, 13).Value = 1.063122641706904
$ws.Cells.Item(17, 14).Value = 1.062567554867882
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.055356267041737
$ws.Cells.Item(18, 4).Value = 1.055680908187163
$ws.Cells.Item(18, 5).Value = 1.059192321931267
$ws.Cells.Item(18, 6).Value = 1.060224039478847
$ws.Cells.Item(18, 9).Value = 1.048084337589997
$ws.Cells.Item(18, 10).Value = 1.061333518655589
$ws.Cells.Item(18, 11).Value = 1.058938545961039
$ws.Cells.Item(18, 12).Value = 1.062438415320103
$ws.Cells.Item(18, 13).Value = 1.063466756849099
$ws.Cells.Item(18, 14).Value = 1.062840733988736
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.055476305357677
$ws.Cells.Item(19, 4).Value = 1.055776679431187
$ws.Cells.Item(19, 5).Value = 1.059308407793279
$ws.Cells.Item(19, 6).Value = 1.060355687208057
$ws.Cells.Item(19, 9).Value = 1.048128288242147
$ws.Cells.Item(19, 10).Value = 1.061426478923031
$ws.Cells.Item(19, 11).Value = 1.059019999251633
$ws.Cells.Item(19, 12).Value = 1.062540171863459
$ws.Cells.Item(19, 13).Value = 1.063584040814661
$ws.Cells.Item(19, 14).Value = 1.062933826270421
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.054939277350045
$ws.Cells.Item(20, 4).Value = 1.055348194853481
$ws.Cells.Item(20, 5).Value = 1.058789085810984
$ws.Cells.Item(20, 6).Value = 1.059766763553434
$ws.Cells.Item(20, 9).Value = 1.047931523940483
$ws.Cells.Item(20, 10).Value = 1.061010522934999
$ws.Cells.Item(20, 11).Value = 1.058655504402964
$ws.Cells.Item(20, 12).Value = 1.062084895208521
$ws.Cells.Item(20, 13).Value = 1.063059320306725
$ws.Cells.Item(20, 14).Value = 1.062517279577175
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.053189990202382
$ws.Cells.Item(21, 4).Value = 1.053952076304611
$ws.Cells.Item(21, 5).Value = 1.057097889324004
$ws.Cells.Item(21, 6).Value = 1.057849181376418
$ws.Cells.Item(21, 9).Value = 1.047288175931405
$ws.Cells.Item(21, 10).Value = 1.059654382546489
$ws.Cells.Item(21, 11).Value = 1.057466660569342
$ws.Cells.Item(21, 12).Value = 1.060601239313482
$ws.Cells.Item(21, 13).Value = 1.061349859265386
$ws.Cells.Item(21, 14).Value = 1.061159213313762
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.052087510345933
$ws.Cells.Item(22, 4).Value = 1.053071881301998
$ws.Cells.Item(22, 5).Value = 1.056032343012097
$ws.Cells.Item(22, 6).Value = 1.056641205306375
$ws.Cells.Item(22, 9).Value = 1.046880856217038
$ws.Cells.Item(22, 10).Value = 1.058798742401892
$ws.Cells.Item(22, 11).Value = 1.056716205829816
$ws.Cells.Item(22, 12).Value = 1.059665667354116
$ws.Cells.Item(22, 13).Value = 1.060272275753269
$ws.Cells.Item(22, 14).Value = 1.060302358062016
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.052672209330508
$ws.Cells.Item(23, 4).Value = 1.053538719652417
$ws.Cells.Item(23, 5).Value = 1.05659742361021
$ws.Cells.Item(23, 6).Value = 1.057281800150762
$ws.Cells.Item(23, 9).Value = 1.04709705314671
$ws.Cells.Item(23, 10).Value = 1.05925261905246
$ws.Cells.Item(23, 11).Value = 1.05711432137197
$ws.Cells.Item(23, 12).Value = 1.060161894596573
$ws.Cells.Item(23, 13).Value = 1.060843790404528
$ws.Cells.Item(23, 14).Value = 1.060756879269488
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.054968560608621
$ws.Cells.Item(24, 4).Value = 1.055371560917526
$ws.Cells.Item(24, 5).Value = 1.058817402030237
$ws.Cells.Item(24, 6).Value = 1.059798873777315
$ws.Cells.Item(24, 9).Value = 1.047942262288613
$ws.Cells.Item(24, 10).Value = 1.061033208960327
$ws.Cells.Item(24, 11).Value = 1.058675385582752
$ws.Cells.Item(24, 12).Value = 1.062109723196357
$ws.Cells.Item(24, 13).Value = 1.063087933446593
$ws.Cells.Item(24, 14).Value = 1.062539997819261
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.057622271520727
$ws.Cells.Item(25, 4).Value = 1.05748830505966
$ws.Cells.Item(25, 5).Value = 1.061384231616193
$ws.Cells.Item(25, 6).Value = 1.062710145417196
$ws.Cells.Item(25, 9).Value = 1.048910956769273
$ws.Cells.Item(25, 10).Value = 1.063086810754482
$ws.Cells.Item(25, 11).Value = 1.060474198607473
$ws.Cells.Item(25, 12).Value = 1.064358469606574
$ws.Cells.Item(25, 13).Value = 1.065680437858125
$ws.Cells.Item(25, 14).Value = 1.064596515963515
